$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.609.40"
$ws.Range("E2").Value = "  +3.72%  "

$ws.Range("D3").Value = "1.913.15"
$ws.Range("E3").Value = "  +1.80%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'248.78"
$ws.Range("E5").Value = "  +1.53%  "

$ws.Range("D6").Value = "'0.698"
$ws.Range("E6").Value = "  +2.86%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'44.26"
$ws.Range("E8").Value = "  +1.06%  "

$ws.Range("D9").Value = "'0.370"
$ws.Range("E9").Value = "  +3.01%  "

$ws.Range("D10").Value = "'57.92"
$ws.Range("E10").Value = "  +8.33%  "

$ws.Range("E11").Value = "  +2.96%  "

$ws.Range("E12").Value = "  +2.54%  "

$ws.Range("D13").Value = "'14.46"
$ws.Range("E13").Value = "  +6.40%  "

$ws.Range("E14").Value = "  +5.44%  "

$ws.Range("D15").Value = "2.191.39"
$ws.Range("E15").Value = "  +1.82%  "

$ws.Range("E16").Value = "  +3.58%  "

$ws.Range("D17").Value = "1.913.24"
$ws.Range("E17").Value = "  +2.19%  "

$ws.Range("D18").Value = "36.616.38"
$ws.Range("E18").Value = "  +3.70%  "

$ws.Range("D19").Value = "'74.40"
$ws.Range("E19").Value = "  +1.41%  "

$ws.Range("E20").Value = "  +4.79%  "

$ws.Range("D21").Value = "'249.91"
$ws.Range("E21").Value = "  +2.35%  "

$ws.Range("E22").Value = "  +4.07%  "

$ws.Range("D23").Value = "'5.19"
$ws.Range("E23").Value = "  +2.96%  "

$ws.Range("D24").Value = "'2.62"
$ws.Range("E24").Value = "  -2.31%  "

$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").Value = "'2.19"
$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("D27").Value = "'168.68"
$ws.Range("E27").Value = "  +2.20%  "

$ws.Range("D28").Value = "'8.81"
$ws.Range("E28").Value = "  +2.51%  "

$ws.Range("D29").Value = "'18.68"
$ws.Range("E29").Value = "  +2.23%  "

$ws.Range("E30").Value = "  +1.47%  "

$ws.Range("D31").Value = "'4.56"
$ws.Range("E31").Value = "  +6.58%  "

$ws.Range("E32").Value = "  +4.69%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.33"
$ws.Range("E33").Value = "  +3.94%  "

$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'1.93"
$ws.Range("E34").Value = "  +6.24%  "

$ws.Range("D35").Value = "'0.0898"
$ws.Range("E35").Value = "  +23.24%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").Value = "'1.51"
$ws.Range("E37").Value = "  +6.69%  "

$ws.Range("D38").Value = "'0.877"
$ws.Range("E38").Value = "  +3.04%  "

$ws.Range("D39").Value = "'17.90"
$ws.Range("E39").Value = "  +53.10%  "

$ws.Range("D40").Value = "'2.03"
$ws.Range("E40").Value = "  +4.52%  "

$ws.Range("D41").Value = "'106.55"
$ws.Range("E41").Value = "  +10.53%  "

$ws.Range("D42").Value = "'0.0226"
$ws.Range("E42").Value = "  +4.28%  "

$ws.Range("D43").Value = "'17.40"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").Value = "'2.95"
$ws.Range("E44").Value = "  +23.24%  "

$ws.Range("D45").Value = "'1.10"
$ws.Range("E45").Value = "  +2.67%  "

$ws.Range("D46").Value = "1.346.37"
$ws.Range("E46").Value = "  +3.07%  "

$ws.Range("D47").Value = "'2.38"
$ws.Range("E47").Value = "  -0.28%  "

$ws.Range("D48").Value = "'0.0812"
$ws.Range("E48").Value = "  +1.83%  "

$ws.Range("D49").Value = "'2.78"
$ws.Range("E49").Value = "  +1.90%  "

$ws.Range("D50").Value = "'6.38"
$ws.Range("E50").Value = "  +1.63%  "

$ws.Range("D51").Value = "'43.39"
$ws.Range("E51").Value = "  +3.13%  "
